$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 461, shifting existing rows 461:544 down to 462:545
$ws.Rows.Item(461).Insert()

# Populate the newly inserted row 461 with its data
$ws.Range("A461").Value = 9
$ws.Range("B461").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C461").Value = "Metropolitana"
$ws.Range("D461").Value = 44995
$ws.Range("E461").Value = 13
$ws.Range("F461").Value = 100112032
$ws.Range("G461").Value = "Zapallo italiano"
$ws.Range("H461").Value = "Sin especificar"
$ws.Range("I461").Value = "Primera"
$ws.Range("J461").Value = 430
$ws.Range("K461").Value = 6000
$ws.Range("L461").Value = 7000
$ws.Range("M461").Value = 6500
$ws.Range("N461").Value = "`$/caja 50 unidades"
$ws.Range("O461").Value = "Región de O'Higgins"
$ws.Range("P461").Value = 130
$ws.Range("Q461").Value = 50
$ws.Range("R461").Value = "Hortaliza"
